# Update instrument reference designators from GA05MOAS-GL002 to GA05MOAS-GL493
# across the Moorings and Asset_Cal_Info sheets, and update the remembered
# selection on the Moorings sheet.

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsCalInfo  = $wb.Worksheets.Item("Asset_Cal_Info")

# Moorings sheet: Ref Des in A2
$wsMoorings.Range("A2").Value = "GA05MOAS-GL493"

# Asset_Cal_Info sheet: Ref Des values in A2:A5, A7, A9, A11
$wsCalInfo.Range("A2").Value = "GA05MOAS-GL493-01-FLORDM000"
$wsCalInfo.Range("A3").Value = "GA05MOAS-GL493-01-FLORDM000"
$wsCalInfo.Range("A4").Value = "GA05MOAS-GL493-01-FLORDM000"
$wsCalInfo.Range("A5").Value = "GA05MOAS-GL493-01-FLORDM000"
$wsCalInfo.Range("A7").Value = "GA05MOAS-GL493-02-DOSTAM000"
$wsCalInfo.Range("A9").Value = "GA05MOAS-GL493-04-CTDGVM000"
$wsCalInfo.Range("A11").Value = "GA05MOAS-GL493-00-ENG000000"

# Update the active selection remembered on the Moorings sheet
$wsMoorings.Activate()
$wsMoorings.Range("F28").Select()
